# Extend the Ground Water init file's data block references from row 18 to
# row 19 (new scenario 19), matching "modified GW init file to extend to
# new scenario 19".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Lower-right cell references for the index blocks: A18/B18/C18/G18/H18/I18/J18
# all need to become A19/B19/C19/G19/H19/I19/J19.
$ws.Range("D5").Value  = "A19"
$ws.Range("D6").Value  = "B19"
$ws.Range("D7").Value  = "C19"
$ws.Range("D8").Value  = "G19"
$ws.Range("D9").Value  = "H19"
$ws.Range("D10").Value = "I19"
$ws.Range("D11").Value = "J19"

# Reflect the author's final selection in the saved view.
$ws.Range("D11").Select()
